$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "JIRAURL"
$ws.Range("B1").Value = "JIRAIssue"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Time"

$ws.Range("B2").Value = "TEST-1"
$ws.Range("C2").Value = (Get-Date -Year 2019 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Value = "3h"

$ws.Hyperlinks.Add($ws.Range("A2"), "http://jira.example.com", "", "", "http://jira.example.com")

$ws.Columns.Item(1).ColumnWidth = 27.833333333333332
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334
